$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new question text in A5
$ws.Range("A5").Value = "How much do you need the delay in the animations to be?"

# Grow the table (Table1) to include the new row
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:B5"))

# Apply formatting to the new row (centered, wrapped, 12pt font) to match the rest of the table
$rng = $ws.Range("A5:B5")
$rng.Font.Size = 12
$rng.HorizontalAlignment = -4108
$rng.WrapText = $true

# Update the active selection to A5, as in the saved workbook
$ws.Range("A5").Select() | Out-Null
